# Adds a new "2022-Q4" sheet (inserted right after "总计", before "2022-Q2")
# and updates the "总计" summary sheet with the new quarter's aggregate row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet before "2022-Q2"
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q4"

# Re-fetch fresh references (the collection shifted after Add()).
$q2 = $wb.Worksheets.Item("2022-Q2")
$q4 = $wb.Worksheets.Item("2022-Q4")

# Copy the header row + index-column formatting from the "2022-Q2" sheet so
# the new sheet matches the look of the other per-quarter sheets.
$q2.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q2.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)

# Headers
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data rows -- keep the code/ratio columns as text (leading zeros, fixed
# decimals) by formatting as Text before writing the values.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "010797"
$q4.Range("C2").Value = "长城优选回报六个月持有期混合A"
$q4.Range("D2").Value = "2.72"
$q4.Range("E2").Value = "31.00"
$q4.Range("F2").Value = "1.32"
$q4.Range("G2").Value = "0.0359"
$q4.Range("H2").Value = 3

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "010798"
$q4.Range("C3").Value = "长城优选回报六个月持有期混合C"
$q4.Range("D3").Value = "0.41"
$q4.Range("E3").Value = "31.00"
$q4.Range("F3").Value = "1.32"
$q4.Range("G3").Value = "0.0054"
$q4.Range("H3").Value = 3

# ---------------------------------------------------------------------------
# 2. Update the "总计" (total) sheet: add the 2022-Q4 row on top and shift
#    the existing rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Add row 6 (2021-Q3, shifted down from row 5), copying the index-column
# style from the row above it.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.02

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.21

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.17

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.6

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04
